# Weekly data refresh for the Economic Dashboard (Aguilar Prototype sheet).
# Updates FRED-sourced present/lag values + 'Latest Date' columns, and promotes
# cells whose date advanced this pull to the yellow 'updated this week' highlight
# (same visual cue already used on N29/N30/N39/N47-N52 etc. in the template).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N3").Value = 45992
$ws.Range("N3").Interior.Color = 65535
$ws.Range("Q3").Value = 50
$ws.Range("R3").Value = 56
$ws.Range("S3").Value = -173
$ws.Range("T3").Value = 108
$ws.Range("U3").Value = -26
$ws.Range("N4").Value = 45992
$ws.Range("N4").Interior.Color = 65535
$ws.Range("Q4").Value = 0.003674296284179135
$ws.Range("R4").Value = 0.005402883639412681
$ws.Range("S4").Value = 0.006706323646421399
$ws.Range("T4").Value = 0.008078881210758367
$ws.Range("U4").Value = 0.008926199121930236
$ws.Range("N6").Value = 45992
$ws.Range("N6").Interior.Color = 65535
$ws.Range("Q6").Value = 4.4
$ws.Range("R6").Value = 4.5
$ws.Range("S6").ClearContents()
$ws.Range("T6").Value = 4.4
$ws.Range("U6").Value = 4.3
$ws.Range("N7").Value = 45992
$ws.Range("N7").Interior.Color = 65535
$ws.Range("Q7").Value = 8.4
$ws.Range("R7").Value = 8.699999999999999
$ws.Range("S7").ClearContents()
$ws.Range("T7").Value = 8.1
$ws.Range("U7").Value = 8.1
$ws.Range("N8").Value = 45992
$ws.Range("N8").Interior.Color = 65535
$ws.Range("Q8").Value = 62.4
$ws.Range("R8").Value = 62.5
$ws.Range("S8").ClearContents()
$ws.Range("T8").Value = 62.5
$ws.Range("U8").Value = 62.3
$ws.Range("N9").Value = 45992
$ws.Range("N9").Interior.Color = 65535
$ws.Range("Q9").Value = 59.7
$ws.Range("R9").Value = 59.6
$ws.Range("S9").ClearContents()
$ws.Range("T9").Value = 59.7
$ws.Range("U9").Value = 59.6
$ws.Range("N15").Value = 45992
$ws.Range("N15").Interior.Color = 65535
$ws.Range("Q15").Value = 34.2
$ws.Range("R15").Value = 34.3
$ws.Range("S15").Value = 34.2
$ws.Range("T15").Value = 34.2
$ws.Range("U15").Value = 34.2
$ws.Range("C24").Value = 45962
$ws.Range("C24").Interior.Color = 65535
$ws.Range("F24").Value = -0.001551215487983981
$ws.Range("G24").Value = 0.004120302067736503
$ws.Range("H24").Value = 0.003065216912115698
$ws.Range("I24").Value = -0.003484844670176512
$ws.Range("J24").Value = 0.008311712980978658
$ws.Range("C25").Value = 45962
$ws.Range("C25").Interior.Color = 65535
$ws.Range("F25").Value = 0.001665686369376251
$ws.Range("G25").Value = 0.001019986599621658
$ws.Range("H25").Value = 0.001988823510608029
$ws.Range("I25").Value = 0.001936089671301433
$ws.Range("J25").Value = 0.0009927073679749654
$ws.Range("N29").Value = 46030
$ws.Range("Q29").Value = 2.23
$ws.Range("R29").Value = 2.24
$ws.Range("S29").Value = 2.24
$ws.Range("T29").Value = 2.23
$ws.Range("U29").Value = 2.22
$ws.Range("N30").Value = 46030
$ws.Range("Q30").Value = 2.27
$ws.Range("R30").Value = 2.27
$ws.Range("S30").Value = 2.27
$ws.Range("T30").Value = 2.26
$ws.Range("U30").Value = 2.25
$ws.Range("N33").Value = 45992
$ws.Range("N33").Interior.Color = 65535
$ws.Range("Q33").Value = 0.003252032520325354
$ws.Range("R33").Value = 0.002444987775060969
$ws.Range("S33").Value = 0.004365620736698661
$ws.Range("T33").Value = 0.001913613996719521
$ws.Range("U33").Value = 0.004117485588800429
$ws.Range("C36").Value = 45931
$ws.Range("C36").Interior.Color = 65535
$ws.Range("F36").Value = 1246
$ws.Range("G36").Value = 1306
$ws.Range("H36").Value = 1291
$ws.Range("I36").Value = 1420
$ws.Range("J36").Value = 1382
$ws.Range("C37").Value = 45931
$ws.Range("C37").Interior.Color = 65535
$ws.Range("F37").Value = -0.07840236686390532
$ws.Range("G37").Value = -0.03758290346352248
$ws.Range("H37").Value = -0.07189072609633357
$ws.Range("I37").Value = 0.1225296442687747
$ws.Range("J37").Value = 0.04144687264506405
$ws.Range("C38").Value = 45931
$ws.Range("C38").Interior.Color = 65535
$ws.Range("F38").Value = 1412
$ws.Range("G38").Value = 1415
$ws.Range("H38").Value = 1330
$ws.Range("I38").Value = 1362
$ws.Range("J38").Value = 1393
$ws.Range("C39").Value = 45931
$ws.Range("C39").Interior.Color = 65535
$ws.Range("F39").Value = -0.01120448179271709
$ws.Range("G39").Value = -0.01324965132496513
$ws.Range("H39").Value = -0.0989159891598916
$ws.Range("I39").Value = -0.05153203342618384
$ws.Range("J39").Value = -0.04654346338124572
$ws.Range("C46").Value = 45931
$ws.Range("C46").Interior.Color = 65535
$ws.Range("F46").Value = 302015
$ws.Range("G46").Value = 294225
$ws.Range("H46").Value = 284060
$ws.Range("I46").Value = 283923
$ws.Range("J46").Value = 280519
$ws.Range("C47").Value = 45931
$ws.Range("C47").Interior.Color = 65535
$ws.Range("F47").Value = 0.02647633613730993
$ws.Range("G47").Value = 0.03578469337463908
$ws.Range("H47").Value = 0.0004825251916893425
$ws.Range("I47").Value = 0.01213465041583639
$ws.Range("J47").Value = 0.0004065548046574552
$ws.Range("N47").Value = 46029
$ws.Range("C48").Value = 45931
$ws.Range("C48").Interior.Color = 65535
$ws.Range("F48").Value = 331366
$ws.Range("G48").Value = 342363
$ws.Range("H48").Value = 339690
$ws.Range("I48").Value = 358321
$ws.Range("J48").Value = 338704
$ws.Range("N48").Value = 46029
$ws.Range("Q48").Value = 3.47
$ws.Range("R48").Value = 3.47
$ws.Range("S48").Value = 3.46
$ws.Range("T48").Value = 3.47
$ws.Range("U48").Value = 3.47
$ws.Range("C49").Value = 45931
$ws.Range("C49").Interior.Color = 65535
$ws.Range("F49").Value = -0.03212087754809956
$ws.Range("G49").Value = 0.007868939327033475
$ws.Range("H49").Value = -0.05199527797700942
$ws.Range("I49").Value = 0.05791782795597333
$ws.Range("J49").Value = -0.03628608905214581
$ws.Range("N49").Value = 46029
$ws.Range("Q49").Value = 3.7
$ws.Range("R49").Value = 3.72
$ws.Range("S49").Value = 3.71
$ws.Range("T49").Value = 3.74
$ws.Range("U49").Value = 3.73
$ws.Range("C50").Value = 45931
$ws.Range("C50").Interior.Color = 65535
$ws.Range("F50").Value = 29796
$ws.Range("G50").Value = 30169
$ws.Range("H50").Value = 30416
$ws.Range("I50").Value = 28606
$ws.Range("J50").Value = 27839
$ws.Range("N50").Value = 46029
$ws.Range("Q50").Value = 4.15
$ws.Range("R50").Value = 4.18
$ws.Range("S50").Value = 4.17
$ws.Range("T50").Value = 4.19
$ws.Range("U50").Value = 4.18
$ws.Range("C51").Value = 45931
$ws.Range("C51").Interior.Color = 65535
$ws.Range("F51").Value = -0.01236368457688353
$ws.Range("G51").Value = -0.008120725933719042
$ws.Range("H51").Value = 0.06327343913864225
$ws.Range("I51").Value = 0.02755127698552395
$ws.Range("J51").Value = 0.05614780530369123
$ws.Range("N51").Value = 46027
$ws.Range("N51").Interior.Color = 65535
$ws.Range("Q51").Value = 6.16
$ws.Range("R51").Value = 6.15
$ws.Range("S51").Value = 6.18
$ws.Range("T51").Value = 6.21
$ws.Range("U51").Value = 6.22
$ws.Range("N52").Value = 46029
$ws.Range("Q52").Value = 5.88
$ws.Range("R52").Value = 5.92
$ws.Range("S52").Value = 5.92
$ws.Range("T52").Value = 5.93
$ws.Range("U52").Value = 5.9
